$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G (RequiredExperience).
# Some become numeric, some become the literal string "Undefined".
$updates = @{
    7   = "Undefined"
    42  = 15
    60  = "Undefined"
    67  = 1
    84  = 1
    85  = 1
    86  = 4
    87  = 3
    90  = 2
    92  = 6
    97  = 3
    104 = 2
    116 = 4
    117 = 10
    119 = 3
    124 = 5
    130 = 4
    132 = 5
    141 = 2
    144 = 5
    145 = 0
    148 = 3
    154 = 2
    163 = "Undefined"
    165 = 5
    166 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
